# Study table.xlsx — "Two more source summaries"
# Adds two new study rows (8 and 9) with their source-paper titles, fixes two
# shared-string typos ("Analytical only" -> "Analytical", "Absorbtion cooling"
# -> "Absorption cooling"), and updates the sheet view / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixes in the row-2 sub-header ---
$ws.Range("C2").Value = "Analytical"
$ws.Range("M2").Value = "Absorption cooling"

# --- Row 8: was a blank placeholder row (just counter "6" in A8) -> becomes
#     a real study row. Keep A8's existing plain style, just change its text.
$ws.Range("A8").Value = "A review of Organic Rankine cycles (ORCs) for the recovery of low-grade waste heat"

# Mark B8, D8, F8, P8 with the checkmark ("ü", Wingdings, centered+bordered)
# style by copying an existing checkmark cell (D3) onto them.
$ws.Range("D3").Copy($ws.Range("B8"))
$ws.Range("D3").Copy($ws.Range("D8"))
$ws.Range("D3").Copy($ws.Range("F8"))
$ws.Range("D3").Copy($ws.Range("P8"))

# M8 is fully removed (no value, no style) in the target.
$ws.Range("M8").Clear() | Out-Null

# Row 8 grows to the same height used by the other study rows.
$ws.Rows.Item(8).RowHeight = 19.5

# --- Row 9: was a blank placeholder row (counter "7" in A9) -> becomes a
#     real study row, with a distinct (Times New Roman 12pt, no border) font.
$ws.Range("A9").Value = "Application of waste heat powered absorption refrigeration system to the LNG recovery process"
$ws.Range("A9").Borders.LineStyle = -4142   # xlLineStyleNone - drop the plain border
$ws.Range("A9").Font.Name = "Times New Roman"
$ws.Range("A9").Font.Size = 12

# Mark B9, M9, P9 with the checkmark style too.
$ws.Range("D3").Copy($ws.Range("B9"))
$ws.Range("D3").Copy($ws.Range("M9"))
$ws.Range("D3").Copy($ws.Range("P9"))

$ws.Rows.Item(9).RowHeight = 19.5

# --- Sheet view: scroll back to A1 (was frozen to show column B first) and
#     move the active selection to B9 (was P7).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B9").Select() | Out-Null
